# TestLog.xlsx minor-update edit script
# Mirrors the "Update maven dependencies and some minor changes" commit:
#  - realign the merged "environment info" banner (A26:D27) from centered to
#    left-aligned (still vertically centered / wrapped)
#  - narrow columns B and D a bit (and drop column B's "best fit" flag)
#  - move the active selection to B36

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan")

# --- Column widths (B narrower, D narrower) -------------------------------
# Excel's COM ColumnWidth is in character units (~1/6-character granularity);
# these are the closest values that round-trip to the target stored widths.
$ws.Columns.Item(2).ColumnWidth = 65.8333333333333
$ws.Columns.Item(4).ColumnWidth = 89

# --- Realign the A26:D27 banner (merged cell block) -----------------------
# -4131 = xlLeft, -4108 = xlCenter
$banner = $ws.Range("A26:D27")
$banner.HorizontalAlignment = -4131
$banner.VerticalAlignment = -4108
$banner.WrapText = $true

# --- Update the active selection -------------------------------------------
$ws.Range("B36").Select()
